$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.953.43'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '1.674.55'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D5").Value = '214.72'
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("E6").Value = '  +1.43%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '0.251'
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("D11").Value = '0.0887'
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = '1.910.71'
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = '1.677.35'
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").Value = '0.527'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '65.80'
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").Value = '26.967.58'
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("D18").Value = '237.39'
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("E19").Value = '  +3.83%  '
$ws.Range("D20").Value = '0.0₃0732'
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("E23").Value = '  -1.24%  '
$ws.Range("E24").Value = '  -2.27%  '
$ws.Range("D25").Value = '145.63'
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = '7.25'
$ws.Range("E26").Value = '  +1.49%  '
$ws.Range("D27").Value = '15.97'
$ws.Range("E27").Value = '  -0.19%  '
$ws.Range("D28").Value = '0.113'
$ws.Range("E28").Value = '  -1.32%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").Value = '0.0497'
$ws.Range("E30").Value = '  -0.29%  '
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = '1.484.27'
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("E35").Value = '  +3.82%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").Value = '0.585'
$ws.Range("E37").Value = '  +1.47%  '
$ws.Range("E38").Value = '  +1.36%  '
$ws.Range("D39").Value = '0.895'
$ws.Range("E39").Value = '  -0.78%  '
$ws.Range("E40").Value = '  -3.08%  '
$ws.Range("E41").Value = '  +2.08%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("E43").Value = '  +1.90%  '
$ws.Range("D44").Value = '66.94'
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("D45").Value = '1.818.52'
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("D46").Value = '0.778'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("D47").Value = '90.42'
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("E50").Value = '  +0.32%  '
$ws.Range("D51").Value = '7.65'
$ws.Range("E51").Value = '  +0.01%  '
